# Actualizacion Datos Personales 4 nov
# Updates statistics (Blancos, Reprobados, Promedio) on sheets
# "Estadisticos 1P" and "Estadisticos Final" for rows 2 (1AV) and 4 (1CV),
# and updates Reprobados on "Estadisticos 2P" for the same rows.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 3
$ws1.Range("E2").Value = 6
$ws1.Range("H2").Value = 7.5
$ws1.Range("D4").Value = 4
$ws1.Range("E4").Value = 1
$ws1.Range("H4").Value = 8.5

# --- Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 37
$ws2.Range("E4").Value = 39

# --- Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 3
$ws3.Range("E2").Value = 6
$ws3.Range("H2").Value = 7.5
$ws3.Range("D4").Value = 4
$ws3.Range("E4").Value = 1
$ws3.Range("H4").Value = 8.5
